$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table to row 45 by cloning row 44 (values + formatting); values overwritten below.
$ws.Range("A44:K44").Copy($ws.Range("A45:K45"))

# Column A holds date-like text (e.g. "16-JAN-26"); force text format so Excel does not
# coerce these into date serial numbers.
$ws.Range("A11:A45").NumberFormat = "@"

# Row 11: 16-JAN-26 - flynas XY-894
$ws.Range("A11").Value = '16-JAN-26'
$ws.Range("B11").Value = 'SM-434'
$ws.Range("C11").Value = 'flynas XY-894'
$ws.Range("D11").Value = 529
$ws.Range("E11").Value = 446
$ws.Range("F11").Value = 83
$ws.Range("G11").Value = 40
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = -10
$ws.Range("J11").Value = 'LOW THREAT'
$ws.Range("K11").Value = 'SAR'

# Row 12: 17-JAN-26 - Nile Air NP-134
$ws.Range("A12").Value = '17-JAN-26'
$ws.Range("B12").Value = 'SM-434'
$ws.Range("C12").Value = 'Nile Air NP-134'
$ws.Range("D12").Value = 368
$ws.Range("E12").Value = 446
$ws.Range("F12").Value = -78
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 'LOW THREAT'
$ws.Range("K12").Value = 'SAR'

# Row 13: 17-JAN-26 - Nile Air NP-132
$ws.Range("A13").Value = '17-JAN-26'
$ws.Range("B13").Value = 'SM-434'
$ws.Range("C13").Value = 'Nile Air NP-132'
$ws.Range("D13").Value = 368
$ws.Range("E13").Value = 446
$ws.Range("F13").Value = -78
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 'LOW THREAT'
$ws.Range("K13").Value = 'SAR'

# Row 14: 17-JAN-26 - Nesma Airlines NE-155
$ws.Range("A14").Value = '17-JAN-26'
$ws.Range("B14").Value = 'SM-434'
$ws.Range("C14").Value = 'Nesma Airlines NE-155'
$ws.Range("D14").Value = 371
$ws.Range("E14").Value = 446
$ws.Range("F14").Value = -75
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 'LOW THREAT'
$ws.Range("K14").Value = 'SAR'

# Row 15: 17-JAN-26 - Nesma Airlines NE-153
$ws.Range("A15").Value = '17-JAN-26'
$ws.Range("B15").Value = 'SM-434'
$ws.Range("C15").Value = 'Nesma Airlines NE-153'
$ws.Range("D15").Value = 371
$ws.Range("E15").Value = 446
$ws.Range("F15").Value = -75
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = 30
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 'LOW THREAT'
$ws.Range("K15").Value = 'SAR'

# Row 16: 17-JAN-26 - Air Arabia Egypt E5-316
$ws.Range("A16").Value = '17-JAN-26'
$ws.Range("B16").Value = 'SM-434'
$ws.Range("C16").Value = 'Air Arabia Egypt E5-316'
$ws.Range("D16").Value = 436
$ws.Range("E16").Value = 446
$ws.Range("F16").Value = -10
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 'LOW THREAT'
$ws.Range("K16").Value = 'SAR'

# Row 17: 17-JAN-26 - flynas XY-854
$ws.Range("A17").Value = '17-JAN-26'
$ws.Range("B17").Value = 'SM-434'
$ws.Range("C17").Value = 'flynas XY-854'
$ws.Range("D17").Value = 489
$ws.Range("E17").Value = 446
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = -10
$ws.Range("J17").Value = 'LOW THREAT'
$ws.Range("K17").Value = 'SAR'

# Row 18: 17-JAN-26 - flynas XY-894
$ws.Range("A18").Value = '17-JAN-26'
$ws.Range("B18").Value = 'SM-434'
$ws.Range("C18").Value = 'flynas XY-894'
$ws.Range("D18").Value = 509
$ws.Range("E18").Value = 446
$ws.Range("F18").Value = 63
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = -10
$ws.Range("J18").Value = 'LOW THREAT'
$ws.Range("K18").Value = 'SAR'

# Row 19: 17-JAN-26 - EgyptAir MS-682
$ws.Range("A19").Value = '17-JAN-26'
$ws.Range("B19").Value = 'SM-434'
$ws.Range("C19").Value = 'EgyptAir MS-682'
$ws.Range("D19").Value = 638
$ws.Range("E19").Value = 446
$ws.Range("F19").Value = 192
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 30
$ws.Range("I19").Value = -16
$ws.Range("J19").Value = 'LOW THREAT'
$ws.Range("K19").Value = 'SAR'

# Row 20: 20-JAN-26 - Nile Air NP-132
$ws.Range("A20").Value = '20-JAN-26'
$ws.Range("B20").Value = 'SM-434'
$ws.Range("C20").Value = 'Nile Air NP-132'
$ws.Range("D20").Value = 368
$ws.Range("E20").Value = 446
$ws.Range("F20").Value = -78
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 'LOW THREAT'
$ws.Range("K20").Value = 'SAR'

# Row 21: 20-JAN-26 - Nesma Airlines NE-151
$ws.Range("A21").Value = '20-JAN-26'
$ws.Range("B21").Value = 'SM-434'
$ws.Range("C21").Value = 'Nesma Airlines NE-151'
$ws.Range("D21").Value = 425
$ws.Range("E21").Value = 446
$ws.Range("F21").Value = -21
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 'LOW THREAT'
$ws.Range("K21").Value = 'SAR'

# Row 22: 20-JAN-26 - flynas XY-854
$ws.Range("A22").Value = '20-JAN-26'
$ws.Range("B22").Value = 'SM-434'
$ws.Range("C22").Value = 'flynas XY-854'
$ws.Range("D22").Value = 489
$ws.Range("E22").Value = 446
$ws.Range("F22").Value = 43
$ws.Range("G22").Value = 40
$ws.Range("H22").Value = 30
$ws.Range("I22").Value = -10
$ws.Range("J22").Value = 'LOW THREAT'
$ws.Range("K22").Value = 'SAR'

# Row 23: 20-JAN-26 - flynas XY-894
$ws.Range("A23").Value = '20-JAN-26'
$ws.Range("B23").Value = 'SM-434'
$ws.Range("C23").Value = 'flynas XY-894'
$ws.Range("D23").Value = 509
$ws.Range("E23").Value = 446
$ws.Range("F23").Value = 63
$ws.Range("G23").Value = 40
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = -10
$ws.Range("J23").Value = 'LOW THREAT'
$ws.Range("K23").Value = 'SAR'

# Row 24: 22-JAN-26 - Nile Air NP-132
$ws.Range("A24").Value = '22-JAN-26'
$ws.Range("B24").Value = 'SM-434'
$ws.Range("C24").Value = 'Nile Air NP-132'
$ws.Range("D24").Value = 368
$ws.Range("E24").Value = 506
$ws.Range("F24").Value = -138
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 'LOW THREAT'
$ws.Range("K24").Value = 'SAR'

# Row 25: 22-JAN-26 - Nesma Airlines NE-155
$ws.Range("A25").Value = '22-JAN-26'
$ws.Range("B25").Value = 'SM-434'
$ws.Range("C25").Value = 'Nesma Airlines NE-155'
$ws.Range("D25").Value = 500
$ws.Range("E25").Value = 506
$ws.Range("F25").Value = -6
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 30
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 'LOW THREAT'
$ws.Range("K25").Value = 'SAR'

# Row 26: 22-JAN-26 - flynas XY-854
$ws.Range("A26").Value = '22-JAN-26'
$ws.Range("B26").Value = 'SM-434'
$ws.Range("C26").Value = 'flynas XY-854'
$ws.Range("D26").Value = 509
$ws.Range("E26").Value = 506
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = -10
$ws.Range("J26").Value = 'LOW THREAT'
$ws.Range("K26").Value = 'SAR'

# Row 27: 23-JAN-26 - Nile Air NP-232
$ws.Range("A27").Value = '23-JAN-26'
$ws.Range("B27").Value = 'SM-434'
$ws.Range("C27").Value = 'Nile Air NP-232'
$ws.Range("D27").Value = 368
$ws.Range("E27").Value = 446
$ws.Range("F27").Value = -78
$ws.Range("G27").Value = 30
$ws.Range("H27").Value = 30
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 'LOW THREAT'
$ws.Range("K27").Value = 'SAR'

# Row 28: 23-JAN-26 - Nesma Airlines NE-153
$ws.Range("A28").Value = '23-JAN-26'
$ws.Range("B28").Value = 'SM-434'
$ws.Range("C28").Value = 'Nesma Airlines NE-153'
$ws.Range("D28").Value = 371
$ws.Range("E28").Value = 446
$ws.Range("F28").Value = -75
$ws.Range("G28").Value = 30
$ws.Range("H28").Value = 30
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 'LOW THREAT'
$ws.Range("K28").Value = 'SAR'

# Row 29: 23-JAN-26 - flynas XY-894
$ws.Range("A29").Value = '23-JAN-26'
$ws.Range("B29").Value = 'SM-434'
$ws.Range("C29").Value = 'flynas XY-894'
$ws.Range("D29").Value = 509
$ws.Range("E29").Value = 446
$ws.Range("F29").Value = 63
$ws.Range("G29").Value = 40
$ws.Range("H29").Value = 30
$ws.Range("I29").Value = -10
$ws.Range("J29").Value = 'LOW THREAT'
$ws.Range("K29").Value = 'SAR'

# Row 30: 24-JAN-26 - Nile Air NP-132
$ws.Range("A30").Value = '24-JAN-26'
$ws.Range("B30").Value = 'SM-434'
$ws.Range("C30").Value = 'Nile Air NP-132'
$ws.Range("D30").Value = 368
$ws.Range("E30").Value = 446
$ws.Range("F30").Value = -78
$ws.Range("G30").Value = 30
$ws.Range("H30").Value = 30
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 'LOW THREAT'
$ws.Range("K30").Value = 'SAR'

# Row 31: 24-JAN-26 - Air Arabia Egypt E5-316
$ws.Range("A31").Value = '24-JAN-26'
$ws.Range("B31").Value = 'SM-434'
$ws.Range("C31").Value = 'Air Arabia Egypt E5-316'
$ws.Range("D31").Value = 436
$ws.Range("E31").Value = 446
$ws.Range("F31").Value = -10
$ws.Range("G31").Value = 30
$ws.Range("H31").Value = 30
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 'LOW THREAT'
$ws.Range("K31").Value = 'SAR'

# Row 32: 24-JAN-26 - flynas XY-854
$ws.Range("A32").Value = '24-JAN-26'
$ws.Range("B32").Value = 'SM-434'
$ws.Range("C32").Value = 'flynas XY-854'
$ws.Range("D32").Value = 509
$ws.Range("E32").Value = 446
$ws.Range("F32").Value = 63
$ws.Range("G32").Value = 40
$ws.Range("H32").Value = 30
$ws.Range("I32").Value = -10
$ws.Range("J32").Value = 'LOW THREAT'
$ws.Range("K32").Value = 'SAR'

# Row 33: 24-JAN-26 - flynas XY-894
$ws.Range("A33").Value = '24-JAN-26'
$ws.Range("B33").Value = 'SM-434'
$ws.Range("C33").Value = 'flynas XY-894'
$ws.Range("D33").Value = 569
$ws.Range("E33").Value = 446
$ws.Range("F33").Value = 123
$ws.Range("G33").Value = 40
$ws.Range("H33").Value = 30
$ws.Range("I33").Value = -10
$ws.Range("J33").Value = 'LOW THREAT'
$ws.Range("K33").Value = 'SAR'

# Row 34: 27-JAN-26 - flynas XY-894
$ws.Range("A34").Value = '27-JAN-26'
$ws.Range("B34").Value = 'SM-434'
$ws.Range("C34").Value = 'flynas XY-894'
$ws.Range("D34").Value = 319
$ws.Range("E34").Value = 476
$ws.Range("F34").Value = -157
$ws.Range("G34").Value = 20
$ws.Range("H34").Value = 30
$ws.Range("I34").Value = 10
$ws.Range("J34").Value = 'LOW THREAT'
$ws.Range("K34").Value = 'SAR'

# Row 35: 27-JAN-26 - flynas XY-854
$ws.Range("A35").Value = '27-JAN-26'
$ws.Range("B35").Value = 'SM-434'
$ws.Range("C35").Value = 'flynas XY-854'
$ws.Range("D35").Value = 319
$ws.Range("E35").Value = 476
$ws.Range("F35").Value = -157
$ws.Range("G35").Value = 20
$ws.Range("H35").Value = 30
$ws.Range("I35").Value = 10
$ws.Range("J35").Value = 'LOW THREAT'
$ws.Range("K35").Value = 'SAR'

# Row 36: 27-JAN-26 - Nile Air NP-132
$ws.Range("A36").Value = '27-JAN-26'
$ws.Range("B36").Value = 'SM-434'
$ws.Range("C36").Value = 'Nile Air NP-132'
$ws.Range("D36").Value = 418
$ws.Range("E36").Value = 476
$ws.Range("F36").Value = -58
$ws.Range("G36").Value = 30
$ws.Range("H36").Value = 30
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 'LOW THREAT'
$ws.Range("K36").Value = 'SAR'

# Row 37: 29-JAN-26 - flynas XY-854
$ws.Range("A37").Value = '29-JAN-26'
$ws.Range("B37").Value = 'SM-434'
$ws.Range("C37").Value = 'flynas XY-854'
$ws.Range("D37").Value = 319
$ws.Range("E37").Value = 476
$ws.Range("F37").Value = -157
$ws.Range("G37").Value = 20
$ws.Range("H37").Value = 30
$ws.Range("I37").Value = 10
$ws.Range("J37").Value = 'LOW THREAT'
$ws.Range("K37").Value = 'SAR'

# Row 38: 29-JAN-26 - Nile Air NP-132
$ws.Range("A38").Value = '29-JAN-26'
$ws.Range("B38").Value = 'SM-434'
$ws.Range("C38").Value = 'Nile Air NP-132'
$ws.Range("D38").Value = 418
$ws.Range("E38").Value = 476
$ws.Range("F38").Value = -58
$ws.Range("G38").Value = 30
$ws.Range("H38").Value = 30
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 'LOW THREAT'
$ws.Range("K38").Value = 'SAR'

# Row 39: 30-JAN-26 - flynas XY-894
$ws.Range("A39").Value = '30-JAN-26'
$ws.Range("B39").Value = 'SM-434'
$ws.Range("C39").Value = 'flynas XY-894'
$ws.Range("D39").Value = 339
$ws.Range("E39").Value = 506
$ws.Range("F39").Value = -167
$ws.Range("G39").Value = 20
$ws.Range("H39").Value = 30
$ws.Range("I39").Value = 10
$ws.Range("J39").Value = 'LOW THREAT'
$ws.Range("K39").Value = 'SAR'

# Row 40: 30-JAN-26 - Nile Air NP-232
$ws.Range("A40").Value = '30-JAN-26'
$ws.Range("B40").Value = 'SM-434'
$ws.Range("C40").Value = 'Nile Air NP-232'
$ws.Range("D40").Value = 418
$ws.Range("E40").Value = 506
$ws.Range("F40").Value = -88
$ws.Range("G40").Value = 30
$ws.Range("H40").Value = 30
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 'LOW THREAT'
$ws.Range("K40").Value = 'SAR'

# Row 41: 30-JAN-26 - Nesma Airlines NE-151
$ws.Range("A41").Value = '30-JAN-26'
$ws.Range("B41").Value = 'SM-434'
$ws.Range("C41").Value = 'Nesma Airlines NE-151'
$ws.Range("D41").Value = 500
$ws.Range("E41").Value = 506
$ws.Range("F41").Value = -6
$ws.Range("G41").Value = 30
$ws.Range("H41").Value = 30
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 'LOW THREAT'
$ws.Range("K41").Value = 'SAR'

# Row 42: 31-JAN-26 - flynas XY-854
$ws.Range("A42").Value = '31-JAN-26'
$ws.Range("B42").Value = 'SM-434'
$ws.Range("C42").Value = 'flynas XY-854'
$ws.Range("D42").Value = 319
$ws.Range("E42").Value = 566
$ws.Range("F42").Value = -247
$ws.Range("G42").Value = 20
$ws.Range("H42").Value = 30
$ws.Range("I42").Value = 10
$ws.Range("J42").Value = 'LOW THREAT'
$ws.Range("K42").Value = 'SAR'

# Row 43: 31-JAN-26 - flynas XY-894
$ws.Range("A43").Value = '31-JAN-26'
$ws.Range("B43").Value = 'SM-434'
$ws.Range("C43").Value = 'flynas XY-894'
$ws.Range("D43").Value = 379
$ws.Range("E43").Value = 566
$ws.Range("F43").Value = -187
$ws.Range("G43").Value = 20
$ws.Range("H43").Value = 30
$ws.Range("I43").Value = 10
$ws.Range("J43").Value = 'LOW THREAT'
$ws.Range("K43").Value = 'SAR'

# Row 44: 31-JAN-26 - Nile Air NP-132
$ws.Range("A44").Value = '31-JAN-26'
$ws.Range("B44").Value = 'SM-434'
$ws.Range("C44").Value = 'Nile Air NP-132'
$ws.Range("D44").Value = 418
$ws.Range("E44").Value = 566
$ws.Range("F44").Value = -148
$ws.Range("G44").Value = 30
$ws.Range("H44").Value = 30
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 'LOW THREAT'
$ws.Range("K44").Value = 'SAR'

# Row 45: 31-JAN-26 - Air Arabia Egypt E5-316
$ws.Range("A45").Value = '31-JAN-26'
$ws.Range("B45").Value = 'SM-434'
$ws.Range("C45").Value = 'Air Arabia Egypt E5-316'
$ws.Range("D45").Value = 436
$ws.Range("E45").Value = 566
$ws.Range("F45").Value = -130
$ws.Range("G45").Value = 30
$ws.Range("H45").Value = 30
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 'LOW THREAT'
$ws.Range("K45").Value = 'SAR'
